# Updated cryptos list on Thu Jun  8 17:41:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.501.11"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5207"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3226"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06774"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.64"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7694"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07778"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.58"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.33"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.011"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.92"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007938"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.544.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.091.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.50%  "

$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.430"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.972"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.27"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.183"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.678"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.71"
$ws.Range("D29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.159"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08731"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.105"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04811"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.866"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.73%  "

$ws.Range("E36").Value = "  +3.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.096"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01782"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.192"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4836"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.26"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8958"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.022"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.603"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4161"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05901"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.048"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.90"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1225"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8844"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.51%  "
